$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# O2 currently holds the image filename text "X071.jpg" for the product that has no
# actual image drawn on the sheet; replace it with a "noimage" placeholder marker.
$ws.Range("O2").Value = "noimage"

# Update the selected / active cell on the sheet to match the author's final position.
$ws.Range("L14").Select()
